$d = $word.ActiveDocument

$old = "We have clarified that the other stands had not been burned in the previous 100 years. The exact date of the last instance of fire at those sites is unknown."
$new = "We have clarified that the other stands had not been burned in the previous 100 years. Although the exact date of the last instance of fire is known in some places throughout Mt. Desert Island, this information was not available at the sites we sampled."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
